# Scheduled-runner refresh of computed profit columns (H:N) across all class sheets.
# Values below were recomputed upstream; this script just writes the refreshed
# numbers (and, where a column newly has/loses a value, adds/clears the cell).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 726.55554
$ws.Cells.Item(2, 9).Value = 165
$ws.Cells.Item(2, 10).Value = 887
$ws.Cells.Item(2, 11).Value = 165
$ws.Cells.Item(2, 12).Value = 887
$ws.Cells.Item(2, 13).Value = -52
$ws.Cells.Item(2, 14).Value = -1113
$ws.Cells.Item(8, 8).Value = 20.75
$ws.Cells.Item(8, 9).Value = 20.75
$ws.Cells.Item(8, 11).Value = 62.25
$ws.Cells.Item(8, 13).Value = 76.75
$ws.Cells.Item(19, 8).Value = 1904.1818
$ws.Cells.Item(19, 9).Value = 1695.5
$ws.Cells.Item(19, 10).Value = 1925.05
$ws.Cells.Item(19, 11).Value = 1695.5
$ws.Cells.Item(19, 12).Value = 1925.05
$ws.Cells.Item(19, 13).Value = -1520.5
$ws.Cells.Item(19, 14).Value = -2275.05
$ws.Cells.Item(28, 8).Value = 368.07693
$ws.Cells.Item(28, 9).Value = 368.07693
$ws.Cells.Item(28, 11).Value = 368.07693
$ws.Cells.Item(28, 13).Value = 116.92307
$ws.Cells.Item(29, 8).Value = 161
$ws.Cells.Item(29, 9).Value = 161
$ws.Cells.Item(29, 11).Value = 483
$ws.Cells.Item(29, 13).Value = -202
$ws.Cells.Item(40, 8).Value = 1740.8334
$ws.Cells.Item(40, 9).Value = 1732.2222
$ws.Cells.Item(40, 10).Value = 1749.4445
$ws.Cells.Item(40, 11).Value = 1732.2222
$ws.Cells.Item(40, 12).Value = 1749.4445
$ws.Cells.Item(40, 13).Value = -1557.2222
$ws.Cells.Item(40, 14).Value = -2099.4445
$ws.Cells.Item(62, 8).Value = 5000
$ws.Cells.Item(62, 9).Value = 5000
$ws.Cells.Item(62, 11).Value = 5000
$ws.Cells.Item(62, 13).Value = -4376
$ws.Cells.Item(65, 8).Value = 5000
$ws.Cells.Item(65, 9).Value = 5000
$ws.Cells.Item(65, 11).Value = 25000
$ws.Cells.Item(65, 13).Value = -21880
$ws.Cells.Item(94, 8).Value = 9988.799999999999
$ws.Cells.Item(94, 9).Value = 9988.799999999999
$ws.Cells.Item(94, 11).Value = 9988.799999999999
$ws.Cells.Item(94, 13).Value = -9537.799999999999
$ws.Cells.Item(98, 8).Value = 1938.5
$ws.Cells.Item(98, 9).Value = 1965.4286
$ws.Cells.Item(98, 11).Value = 1965.4286
$ws.Cells.Item(98, 13).Value = -467.4286
$ws.Cells.Item(100, 8).Value = 2147.3333
$ws.Cells.Item(100, 9).Value = 1971
$ws.Cells.Item(100, 10).Value = 2500
$ws.Cells.Item(100, 11).Value = 1971
$ws.Cells.Item(100, 12).Value = 2500
$ws.Cells.Item(100, 13).Value = -1430
$ws.Cells.Item(100, 14).Value = -3582
$ws.Cells.Item(122, 8).Value = 1938.5
$ws.Cells.Item(122, 9).Value = 1965.4286
$ws.Cells.Item(122, 11).Value = 5896.2858
$ws.Cells.Item(122, 13).Value = -3446.2858
$ws.Cells.Item(137, 8).Value = 2152.3333
$ws.Cells.Item(137, 10).Value = 2230
$ws.Cells.Item(137, 12).Value = 6690
$ws.Cells.Item(137, 14).Value = -11790
$ws.Cells.Item(138, 8).Value = 2171.68
$ws.Cells.Item(138, 10).Value = 3327.5715
$ws.Cells.Item(138, 12).Value = 9982.7145
$ws.Cells.Item(138, 14).Value = -20262.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5239.4414
$ws.Cells.Item(32, 9).Value = 3611.8215
$ws.Cells.Item(32, 11).Value = 3611.8215
$ws.Cells.Item(32, 13).Value = -3324.8215
$ws.Cells.Item(45, 8).Value = 2399.3333
$ws.Cells.Item(45, 9).Value = 1798
$ws.Cells.Item(45, 11).Value = 1798
$ws.Cells.Item(45, 13).Value = -1421
$ws.Cells.Item(61, 8).Value = 3934.25
$ws.Cells.Item(61, 9).Value = 3966.1738
$ws.Cells.Item(61, 11).Value = 3966.1738
$ws.Cells.Item(61, 13).Value = -3754.1738
$ws.Cells.Item(74, 8).Value = 858
$ws.Cells.Item(74, 9).Value = 684.3077
$ws.Cells.Item(74, 11).Value = 684.3077
$ws.Cells.Item(74, 13).Value = 189.6923
$ws.Cells.Item(77, 8).Value = 858
$ws.Cells.Item(77, 9).Value = 684.3077
$ws.Cells.Item(77, 11).Value = 3421.5385
$ws.Cells.Item(77, 13).Value = 946.4615000000003
$ws.Cells.Item(110, 8).Value = 2692.182
$ws.Cells.Item(110, 9).Value = 2201.75
$ws.Cells.Item(110, 11).Value = 2201.75
$ws.Cells.Item(110, 13).Value = -156.75
$ws.Cells.Item(122, 8).Value = 538088.8
$ws.Cells.Item(122, 9).Value = 673584.5600000001
$ws.Cells.Item(122, 11).Value = 2020753.68
$ws.Cells.Item(122, 13).Value = -2018303.68
$ws.Cells.Item(136, 8).Value = 3934.25
$ws.Cells.Item(136, 9).Value = 3966.1738
$ws.Cells.Item(136, 11).Value = 11898.5214
$ws.Cells.Item(136, 13).Value = -9348.5214

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 2719.5
$ws.Cells.Item(22, 9).Value = 292.66666
$ws.Cells.Item(22, 11).Value = 292.66666
$ws.Cells.Item(22, 13).Value = -119.66666
$ws.Cells.Item(35, 8).Value = 21174.5
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).ClearContents()
$ws.Cells.Item(107, 8).Value = 1000
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 13).Value = 920

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 113.333336
$ws.Cells.Item(10, 9).Value = 95
$ws.Cells.Item(10, 10).Value = 150
$ws.Cells.Item(10, 11).Value = 95
$ws.Cells.Item(10, 12).Value = 150
$ws.Cells.Item(10, 13).Value = 44
$ws.Cells.Item(10, 14).Value = -428
$ws.Cells.Item(22, 8).Value = 6941
$ws.Cells.Item(22, 9).Value = 237.90909
$ws.Cells.Item(22, 10).Value = 25374.5
$ws.Cells.Item(22, 11).Value = 237.90909
$ws.Cells.Item(22, 12).Value = 25374.5
$ws.Cells.Item(22, 13).Value = 112.09091
$ws.Cells.Item(22, 14).Value = -26074.5
$ws.Cells.Item(31, 8).Value = 5272.3
$ws.Cells.Item(31, 9).Value = 5272.3
$ws.Cells.Item(31, 11).Value = 5272.3
$ws.Cells.Item(31, 13).Value = -4977.3
$ws.Cells.Item(34, 8).Value = 5272.3
$ws.Cells.Item(34, 9).Value = 5272.3
$ws.Cells.Item(34, 11).Value = 5272.3
$ws.Cells.Item(34, 13).Value = -5070.3
$ws.Cells.Item(94, 8).Value = 4250.25
$ws.Cells.Item(94, 10).Value = 4500.3335
$ws.Cells.Item(94, 12).Value = 4500.3335
$ws.Cells.Item(94, 14).Value = -5402.3335
$ws.Cells.Item(99, 8).Value = 11658.23
$ws.Cells.Item(99, 9).Value = 7259.9287
$ws.Cells.Item(99, 10).Value = 16789.584
$ws.Cells.Item(99, 11).Value = 7259.9287
$ws.Cells.Item(99, 12).Value = 16789.584
$ws.Cells.Item(99, 13).Value = -5761.9287
$ws.Cells.Item(99, 14).Value = -19785.584
$ws.Cells.Item(126, 8).Value = 11658.23
$ws.Cells.Item(126, 9).Value = 7259.9287
$ws.Cells.Item(126, 10).Value = 16789.584
$ws.Cells.Item(126, 11).Value = 21779.7861
$ws.Cells.Item(126, 12).Value = 50368.75199999999
$ws.Cells.Item(126, 13).Value = -19309.7861
$ws.Cells.Item(126, 14).Value = -55308.75199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 7744768.5
$ws.Cells.Item(32, 10).Value = 8712614
$ws.Cells.Item(32, 12).Value = 26137842
$ws.Cells.Item(32, 14).Value = -26138408
$ws.Cells.Item(33, 8).Value = 112302.22
$ws.Cells.Item(33, 9).Value = 1753.5
$ws.Cells.Item(33, 10).Value = 333399.66
$ws.Cells.Item(33, 11).Value = 10521
$ws.Cells.Item(33, 12).Value = 2000397.96
$ws.Cells.Item(33, 13).Value = -10238
$ws.Cells.Item(33, 14).Value = -2000963.96
$ws.Cells.Item(38, 8).Value = 76.53846
$ws.Cells.Item(38, 9).Value = 50.545456
$ws.Cells.Item(38, 11).Value = 151.636368
$ws.Cells.Item(38, 13).Value = 195.363632
$ws.Cells.Item(40, 8).Value = 118.6
$ws.Cells.Item(40, 9).Value = 181
$ws.Cells.Item(40, 10).Value = 25
$ws.Cells.Item(40, 11).Value = 724
$ws.Cells.Item(40, 12).Value = 100
$ws.Cells.Item(40, 13).Value = -655
$ws.Cells.Item(40, 14).Value = -238
$ws.Cells.Item(120, 8).Value = 1799.3334
$ws.Cells.Item(120, 10).Value = 5000
$ws.Cells.Item(120, 12).Value = 15000
$ws.Cells.Item(120, 14).Value = -24676
$ws.Cells.Item(131, 8).Value = 1454.8727
$ws.Cells.Item(131, 10).Value = 1519.76
$ws.Cells.Item(131, 12).Value = 4559.28
$ws.Cells.Item(131, 14).Value = -14639.28
$ws.Cells.Item(136, 8).Value = 10499.5
$ws.Cells.Item(136, 9).Value = 999
$ws.Cells.Item(136, 11).Value = 2997
$ws.Cells.Item(136, 13).Value = 2103
$ws.Cells.Item(138, 8).Value = 2210.1538
$ws.Cells.Item(138, 9).Value = 1252.8334
$ws.Cells.Item(138, 11).Value = 3758.5002
$ws.Cells.Item(138, 13).Value = 1381.4998
$ws.Cells.Item(139, 8).Value = 3124.75
$ws.Cells.Item(139, 9).Value = 3333.3333
$ws.Cells.Item(139, 11).Value = 9999.999899999999
$ws.Cells.Item(139, 13).Value = -4859.999899999999
$ws.Cells.Item(140, 8).Value = 3411.4707
$ws.Cells.Item(140, 9).Value = 2545.4546
$ws.Cells.Item(140, 11).Value = 7636.3638
$ws.Cells.Item(140, 13).Value = -2456.3638
$ws.Cells.Item(141, 8).Value = 5008.5713
$ws.Cells.Item(141, 9).Value = 5008.5713
$ws.Cells.Item(141, 11).Value = 15025.7139
$ws.Cells.Item(141, 13).Value = -9845.713899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3909.7
$ws.Cells.Item(80, 9).Value = 1650
$ws.Cells.Item(80, 10).Value = 4474.625
$ws.Cells.Item(80, 11).Value = 1650
$ws.Cells.Item(80, 12).Value = 4474.625
$ws.Cells.Item(80, 13).Value = -652
$ws.Cells.Item(80, 14).Value = -6470.625
$ws.Cells.Item(83, 8).Value = 3909.7
$ws.Cells.Item(83, 9).Value = 1650
$ws.Cells.Item(83, 10).Value = 4474.625
$ws.Cells.Item(83, 11).Value = 8250
$ws.Cells.Item(83, 12).Value = 22373.125
$ws.Cells.Item(83, 13).Value = -3258
$ws.Cells.Item(83, 14).Value = -32357.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5651.4546
$ws.Cells.Item(7, 9).Value = 5966.6
$ws.Cells.Item(7, 11).Value = 5966.6
$ws.Cells.Item(7, 13).Value = -5854.6
$ws.Cells.Item(40, 8).Value = 2250
$ws.Cells.Item(40, 9).Value = 2250
$ws.Cells.Item(40, 11).Value = 2250
$ws.Cells.Item(40, 13).Value = -2114
$ws.Cells.Item(47, 8).Value = 18000
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()
$ws.Cells.Item(52, 8).Value = 18000
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 5651.4546
$ws.Cells.Item(126, 9).Value = 5966.6
$ws.Cells.Item(126, 11).Value = 17899.8
$ws.Cells.Item(126, 13).Value = -15429.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2348.875
$ws.Cells.Item(126, 9).Value = 1949.25
$ws.Cells.Item(126, 10).Value = 2748.5
$ws.Cells.Item(126, 11).Value = 5847.75
$ws.Cells.Item(126, 12).Value = 8245.5
$ws.Cells.Item(126, 13).Value = -3377.75
$ws.Cells.Item(126, 14).Value = -13185.5
